# edit.ps1 - applies the Brasil COVID dataset update (FT plot fix + 4 new days of data)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Confirmados")
$ws2 = $wb.Worksheets.Item("Mortes")

# --- Corrections to existing rows (125, 127, 128) ---
$ws1.Range("N125").Value = 29784
$ws1.Range("N127").Value = 31964
$ws1.Range("J128").Value = 26145
$ws1.Range("N128").Value = 32962
$ws1.Range("Q128").Value = 9062
$ws1.Range("R128").Value = 26612
$ws1.Range("V128").Value = 11454
$ws1.Range("Z128").Value = 17401
$ws1.Range("AA128").Value = 21523
$ws1.Range("AB128").Value = 17583

$ws2.Range("N125").Value = 999
$ws2.Range("N127").Value = 1072
$ws2.Range("N128").Value = 1108
$ws2.Range("Q128").Value = 91
$ws2.Range("R128").Value = 725
$ws2.Range("V128").Value = 209
$ws2.Range("Z128").Value = 665
$ws2.Range("AA128").Value = 523
$ws2.Range("AB128").Value = 354

# --- New rows 129-132 (2020-07-02 .. 2020-07-05) on Confirmados ---
# Row 129: 2020-07-02
$data129 = @(302179,50242,116823,79349,52281,37328,50707,29195,26318,49536,61119,26304,33487,28575,73530,9388,26612,116519,14048,110411,11736,24376,29153,86025,18356,21970,18323)
$vals129 = New-Object 'object[,]' 1,27
for ($i=0; $i -lt 27; $i++) { $vals129[0,$i] = $data129[$i] }
$ws1.Range("B129:AB129").Value = $vals129
$c = $ws1.Range("A129")
$c.Formula = '="2020-07-02"'
$c.Copy()
$c.PasteSpecial(-4163)

# Row 130: 2020-07-03
$data130 = @(310517,51689,118956,82314,53996,38404,53351,30371,27502,50765,62362,28166,33910,30261,74537,9910,28186,118311,14112,112531,12282,25561,29574,88214,19540,22241,18769)
$vals130 = New-Object 'object[,]' 1,27
for ($i=0; $i -lt 27; $i++) { $vals130[0,$i] = $data130[$i] }
$ws1.Range("B130:AB130").Value = $vals130
$c = $ws1.Range("A130")
$c.Formula = '="2020-07-03"'
$c.Copy()
$c.PasteSpecial(-4163)

# Row 131: 2020-07-04
$data131 = @(312530,52824,120428,85485,55760,39255,55958,31619,28526,52306,63457,30570,34645,31931,75945,10089,29761,121464,14487,113811,12475,26079,29809,89057,20333,22957,18922)
$vals131 = New-Object 'object[,]' 1,27
for ($i=0; $i -lt 27; $i++) { $vals131[0,$i] = $data131[$i] }
$ws1.Range("B131:AB131").Value = $vals131
$c = $ws1.Range("A131")
$c.Formula = '="2020-07-04"'
$c.Copy()
$c.PasteSpecial(-4163)

# Row 132: 2020-07-05
$data132 = @(320179,53393,121292,87048,57854,39935,58283,31955,28526,52728,65129,31459,34645,32969,76014,10089,30217,121986,14622,114535,12475,26511,29883,89714,21081,23479,18922)
$vals132 = New-Object 'object[,]' 1,27
for ($i=0; $i -lt 27; $i++) { $vals132[0,$i] = $data132[$i] }
$ws1.Range("B132:AB132").Value = $vals132
$c = $ws1.Range("A132")
$c.Formula = '="2020-07-05"'
$c.Copy()
$c.PasteSpecial(-4163)

# --- New rows 129-132 (2020-07-02 .. 2020-07-05) on Mortes ---
# Row 129: 2020-07-02
$data129 = @(15351,1727,10332,1947,631,1091,1059,663,572,1044,4968,706,1177,362,2862,107,725,6307,378,5050,211,726,427,2119,706,530,358)
$vals129 = New-Object 'object[,]' 1,27
for ($i=0; $i -lt 27; $i++) { $vals129[0,$i] = $data129[$i] }
$ws2.Range("B129:AB129").Value = $vals129
$c = $ws2.Range("A129")
$c.Formula = '="2020-07-02"'
$c.Copy()
$c.PasteSpecial(-4163)

# Row 130: 2020-07-03
$data130 = @(15694,1758,10500,2001,643,1113,1110,690,602,1062,5068,728,1200,376,2887,114,764,6373,387,5069,215,762,438,2153,741,533,364)
$vals130 = New-Object 'object[,]' 1,27
for ($i=0; $i -lt 27; $i++) { $vals130[0,$i] = $data130[$i] }
$ws2.Range("B130:AB130").Value = $vals130
$c = $ws2.Range("A130")
$c.Formula = '="2020-07-03"'
$c.Copy()
$c.PasteSpecial(-4163)

# Row 131: 2020-07-04
$data131 = @(15996,1781,10624,2050,671,1134,1183,715,621,1082,5116,763,1213,383,2918,117,783,6441,391,5096,220,784,441,2185,786,548,368)
$vals131 = New-Object 'object[,]' 1,27
for ($i=0; $i -lt 27; $i++) { $vals131[0,$i] = $data131[$i] }
$ws2.Range("B131:AB131").Value = $vals131
$c = $ws2.Range("A131")
$c.Formula = '="2020-07-04"'
$c.Copy()
$c.PasteSpecial(-4163)

# Row 132: 2020-07-05
$data132 = @(16078,1803,10667,2107,699,1153,1201,727,621,1099,5143,795,1213,393,2929,117,798,6441,391,5105,220,809,442,2219,821,555,368)
$vals132 = New-Object 'object[,]' 1,27
for ($i=0; $i -lt 27; $i++) { $vals132[0,$i] = $data132[$i] }
$ws2.Range("B132:AB132").Value = $vals132
$c = $ws2.Range("A132")
$c.Formula = '="2020-07-05"'
$c.Copy()
$c.PasteSpecial(-4163)

$excel.CutCopyMode = $false

Write-Host "Done applying Brasil dataset update."